$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (copy style/format from the existing header cell AC1
# so the new headers keep the same bold/border/center formatting).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row (2-42).
$wins = 83
$losses = 79
$ties = 0

for ($row = 2; $row -le 42; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins
    $ws.Cells.Item($row, 31).Value = $losses
    $ws.Cells.Item($row, 32).Value = $ties
}
